$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 71, pushing existing rows 71-100 down to 72-101
$ws.Rows("71:71").Insert()

# Populate the newly inserted row 71 with the new weekly record
$ws.Range("A71").Value = 10
$ws.Range("B71").Value = "Vega Modelo de Temuco"
$ws.Range("C71").Value = "La Araucanía"
$ws.Range("D71").Value = 44755
$ws.Range("E71").Value = 9
$ws.Range("F71").Value = 100114002
$ws.Range("G71").Value = "Camote"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 30
$ws.Range("K71").Value = 20000
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = 20000
$ws.Range("N71").Value = "$/malla 20 kilos"
$ws.Range("O71").Value = "Perú"
$ws.Range("P71").Value = 1000
$ws.Range("Q71").Value = 20
$ws.Range("R71").Value = "Hortaliza"
